# Updated Jira ids for STeAM test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Jira id renames (A2 keeps its existing style) ---
$ws.Range("A2").Value = "OPQA-1410"

# --- Column B: reworded description for the "existing user" test case ---
$ws.Range("B3").Value = "Verify that existing user can't be created and check the error status using STeAM API "

# A3, A5, A6, A7, A8, A9, A10 lose their explicit style (reset to Normal) as part of the edit;
# A4 keeps its existing style.
$ws.Range("A3").Style = "Normal"
$ws.Range("A3").Value = "OPQA-1411"

$ws.Range("A4").Value = "OPQA-1412"

$ws.Range("A5").Style = "Normal"
$ws.Range("A5").Value = "OPQA-1413"

$ws.Range("A6").Style = "Normal"
$ws.Range("A6").Value = "OPQA-1414"

$ws.Range("A7").Style = "Normal"
$ws.Range("A7").Value = "OPQA-1415"

$ws.Range("A8").Style = "Normal"
$ws.Range("A8").Value = "OPQA-1416"

$ws.Range("A9").Style = "Normal"
$ws.Range("A9").Value = "OPQA-1417"

$ws.Range("A10").Style = "Normal"
$ws.Range("A10").Value = "OPQA-1418"

# --- Column K: VALIDATIONS strings embedding the renamed ids ---
$ws.Range("K4").Value = "status=200||rc=OK||USER.USER_ID=(OPQA-1410_user.userID)||UserInfo.USER_INFO_FIRST_NAME=JANARDHANUpadateF||UserInfo.USER_INFO_MIDDLE_NAME=REDDYUpdatedM||UserInfo.USER_INFO_LAST_NAME=UpdatedL"
$ws.Range("K7").Value = "status=200||rc=OK||user.userID=(OPQA-1410_user.userID)"
$ws.Range("K8").Value = "status=200||rc=OK||USER.UserId=(OPQA-1410_user.userID)"
$ws.Range("K9").Value = "status=200||rc=OK||Results.USER_ID=(OPQA-1410_user.userID)||Results.TRUID=(OPQA-1410_user.truID)||Results.TOTAL_ROWS=1"
$ws.Range("K6").Value = "status=200||rc=OK||User.userID=(OPQA-1410_user.userID)||User.truId=(OPQA-1410_user.truID)||UserInfo.USER_INFO_FIRST_NAME=(OPQA-1412_UserInfo.USER_INFO_FIRST_NAME)||UserInfo.USER_INFO_MIDDLE_NAME=(OPQA-1412_UserInfo.USER_INFO_MIDDLE_NAME)||UserInfo.USER_INFO_LAST_NAME=(OPQA-1412_UserInfo.USER_INFO_LAST_NAME)"

# --- Column J: DEPENDENCYTESTS references to the renamed OPQA-AAA id (reuses the OPQA-1410 string) ---
$ws.Range("J3").Value = "OPQA-1410"
$ws.Range("J4").Value = "OPQA-1410"
$ws.Range("J5").Value = "OPQA-1410"
$ws.Range("J6").Value = "OPQA-1410"
$ws.Range("J7").Value = "OPQA-1410"
$ws.Range("J8").Value = "OPQA-1410"
$ws.Range("J9").Value = "OPQA-1410"
$ws.Range("J10").Value = "OPQA-1410"

# --- Selection: active cell within M2:M10 moved from M10 to M2 ---
$ws.Range("M2:M10").Select()
